# performance_es_2.xlsx update:
# - add CLOSENESS NAIVE timing result (D3)
# - add CLOSENESS PARALLEL timing results (D11, E11)
# - widen columns A, B and D to fit the new/longer content
# - leave the active cell selection on E20 (as last interacted by the author)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New data points
$ws.Range("D3").Value = 1554.2362387180301
$ws.Range("D11").Value = 269.693302631378
$ws.Range("E11").Value = 265.10288429260203

# Column widths (character units; Excel snaps these to whole pixels using the
# workbook's default font metrics, so the values below are chosen to land as
# close as possible to the widths recorded for columns A, B and D).
$ws.Columns.Item(1).ColumnWidth = 33.67
$ws.Columns.Item(2).ColumnWidth = 25.83
$ws.Columns.Item(4).ColumnWidth = 11.5

# Final selection/active cell
$ws.Range("E20").Select()
